$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Sistema B"
